$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by all rows in this block (market/region/product info)
$constA = 7
$constB = 'Terminal Hortofrutícola Agro Chillán'
$constC = 'Ñuble'
$constE = 16
$constF = 'Fruta'
$constG = 100108
$constH = 'Tropicales y subtropicales'
$constI = 100108005
$constJ = 'Piña'
$constK = 'Caramelo'
$constR = 'Ecuador'

# Populate the two brand-new rows (121, 122) with the shared constant columns
foreach ($r in 121..122) {
    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
    $ws.Cells.Item($r, 10).Value = $constJ
    $ws.Cells.Item($r, 11).Value = $constK
    $ws.Cells.Item($r, 18).Value = $constR
}

# Apply the same date style (s=2 / custom date number format) used by column D in this block to the new rows
$ws.Range("D121").NumberFormat = $ws.Range("D119").NumberFormat
$ws.Range("D122").NumberFormat = $ws.Range("D119").NumberFormat

# Row-by-row data update: Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad, Precio $/Kg, Kg/unidad
$ws.Range("D72").Value = 44447
$ws.Range("L72").Value = 'Primera'
$ws.Range("M72").Value = 60
$ws.Range("N72").Value = 19000
$ws.Range("O72").Value = 20000
$ws.Range("P72").Value = 19500
$ws.Range("Q72").Value = '$/caja 12 unidades'
$ws.Range("S72").Value = 1625
$ws.Range("T72").Value = 12

$ws.Range("D73").Value = 44447
$ws.Range("L73").Value = 'Segunda'
$ws.Range("M73").Value = 60
$ws.Range("N73").Value = 19000
$ws.Range("O73").Value = 20000
$ws.Range("P73").Value = 19500
$ws.Range("Q73").Value = '$/caja 14 unidades'
$ws.Range("S73").Value = 1393
$ws.Range("T73").Value = 14

$ws.Range("D74").Value = 44267
$ws.Range("L74").Value = 'Segunda'
$ws.Range("M74").Value = 120
$ws.Range("N74").Value = 15000
$ws.Range("O74").Value = 16000
$ws.Range("P74").Value = 15500
$ws.Range("Q74").Value = '$/caja 14 unidades'
$ws.Range("S74").Value = 1107
$ws.Range("T74").Value = 14

$ws.Range("D75").Value = 44292
$ws.Range("L75").Value = 'Segunda'
$ws.Range("M75").Value = 60
$ws.Range("N75").Value = 15500
$ws.Range("O75").Value = 16000
$ws.Range("P75").Value = 15750
$ws.Range("Q75").Value = '$/caja 14 unidades'
$ws.Range("S75").Value = 1125
$ws.Range("T75").Value = 14

$ws.Range("D76").Value = 44300
$ws.Range("L76").Value = 'Segunda'
$ws.Range("M76").Value = 120
$ws.Range("N76").Value = 15000
$ws.Range("O76").Value = 16000
$ws.Range("P76").Value = 15500
$ws.Range("Q76").Value = '$/caja 14 unidades'
$ws.Range("S76").Value = 1107
$ws.Range("T76").Value = 14

$ws.Range("D77").Value = 44392
$ws.Range("L77").Value = 'Segunda'
$ws.Range("M77").Value = 120
$ws.Range("N77").Value = 16000
$ws.Range("O77").Value = 17000
$ws.Range("P77").Value = 16500
$ws.Range("Q77").Value = '$/caja 14 unidades'
$ws.Range("S77").Value = 1179
$ws.Range("T77").Value = 14

$ws.Range("D78").Value = 44224
$ws.Range("L78").Value = 'Segunda'
$ws.Range("M78").Value = 30
$ws.Range("N78").Value = 15000
$ws.Range("O78").Value = 16000
$ws.Range("P78").Value = 15500
$ws.Range("Q78").Value = '$/caja 14 unidades'
$ws.Range("S78").Value = 1107
$ws.Range("T78").Value = 14

$ws.Range("D79").Value = 44295
$ws.Range("L79").Value = 'Segunda'
$ws.Range("M79").Value = 60
$ws.Range("N79").Value = 15000
$ws.Range("O79").Value = 16000
$ws.Range("P79").Value = 15500
$ws.Range("Q79").Value = '$/caja 14 unidades'
$ws.Range("S79").Value = 1107
$ws.Range("T79").Value = 14

$ws.Range("D80").Value = 44179
$ws.Range("L80").Value = 'Primera'
$ws.Range("M80").Value = 45
$ws.Range("N80").Value = 17000
$ws.Range("O80").Value = 18000
$ws.Range("P80").Value = 17556
$ws.Range("Q80").Value = '$/caja 12 unidades'
$ws.Range("S80").Value = 1463
$ws.Range("T80").Value = 12

$ws.Range("D81").Value = 44179
$ws.Range("L81").Value = 'Segunda'
$ws.Range("M81").Value = 50
$ws.Range("N81").Value = 17000
$ws.Range("O81").Value = 18000
$ws.Range("P81").Value = 17400
$ws.Range("Q81").Value = '$/caja 14 unidades'
$ws.Range("S81").Value = 1243
$ws.Range("T81").Value = 14

$ws.Range("D82").Value = 44362
$ws.Range("L82").Value = 'Segunda'
$ws.Range("M82").Value = 120
$ws.Range("N82").Value = 15000
$ws.Range("O82").Value = 16000
$ws.Range("P82").Value = 15500
$ws.Range("Q82").Value = '$/caja 14 unidades'
$ws.Range("S82").Value = 1107
$ws.Range("T82").Value = 14

$ws.Range("D83").Value = 44309
$ws.Range("L83").Value = 'Segunda'
$ws.Range("M83").Value = 60
$ws.Range("N83").Value = 15000
$ws.Range("O83").Value = 16000
$ws.Range("P83").Value = 15500
$ws.Range("Q83").Value = '$/caja 14 unidades'
$ws.Range("S83").Value = 1107
$ws.Range("T83").Value = 14

$ws.Range("D84").Value = 44384
$ws.Range("L84").Value = 'Segunda'
$ws.Range("M84").Value = 120
$ws.Range("N84").Value = 16500
$ws.Range("O84").Value = 17000
$ws.Range("P84").Value = 16750
$ws.Range("Q84").Value = '$/caja 14 unidades'
$ws.Range("S84").Value = 1196
$ws.Range("T84").Value = 14

$ws.Range("D85").Value = 44196
$ws.Range("L85").Value = 'Primera'
$ws.Range("M85").Value = 50
$ws.Range("N85").Value = 14000
$ws.Range("O85").Value = 15000
$ws.Range("P85").Value = 14500
$ws.Range("Q85").Value = '$/caja 12 unidades'
$ws.Range("S85").Value = 1208
$ws.Range("T85").Value = 12

$ws.Range("D86").Value = 44196
$ws.Range("L86").Value = 'Segunda'
$ws.Range("M86").Value = 50
$ws.Range("N86").Value = 14000
$ws.Range("O86").Value = 15000
$ws.Range("P86").Value = 14400
$ws.Range("Q86").Value = '$/caja 14 unidades'
$ws.Range("S86").Value = 1029
$ws.Range("T86").Value = 14

$ws.Range("D87").Value = 44365
$ws.Range("L87").Value = 'Segunda'
$ws.Range("M87").Value = 120
$ws.Range("N87").Value = 17000
$ws.Range("O87").Value = 18000
$ws.Range("P87").Value = 17500
$ws.Range("Q87").Value = '$/caja 14 unidades'
$ws.Range("S87").Value = 1250
$ws.Range("T87").Value = 14

$ws.Range("D88").Value = 44218
$ws.Range("L88").Value = 'Primera'
$ws.Range("M88").Value = 55
$ws.Range("N88").Value = 13500
$ws.Range("O88").Value = 14500
$ws.Range("P88").Value = 13955
$ws.Range("Q88").Value = '$/caja 12 unidades'
$ws.Range("S88").Value = 1163
$ws.Range("T88").Value = 12

$ws.Range("D89").Value = 44321
$ws.Range("L89").Value = 'Segunda'
$ws.Range("M89").Value = 120
$ws.Range("N89").Value = 15000
$ws.Range("O89").Value = 16000
$ws.Range("P89").Value = 15500
$ws.Range("Q89").Value = '$/caja 14 unidades'
$ws.Range("S89").Value = 1107
$ws.Range("T89").Value = 14

$ws.Range("D90").Value = 44291
$ws.Range("L90").Value = 'Segunda'
$ws.Range("M90").Value = 120
$ws.Range("N90").Value = 15000
$ws.Range("O90").Value = 16000
$ws.Range("P90").Value = 15500
$ws.Range("Q90").Value = '$/caja 14 unidades'
$ws.Range("S90").Value = 1107
$ws.Range("T90").Value = 14

$ws.Range("D91").Value = 44215
$ws.Range("L91").Value = 'Primera'
$ws.Range("M91").Value = 65
$ws.Range("N91").Value = 12500
$ws.Range("O91").Value = 13000
$ws.Range("P91").Value = 12731
$ws.Range("Q91").Value = '$/caja 12 unidades'
$ws.Range("S91").Value = 1061
$ws.Range("T91").Value = 12

$ws.Range("D92").Value = 44194
$ws.Range("L92").Value = 'Segunda'
$ws.Range("M92").Value = 120
$ws.Range("N92").Value = 14000
$ws.Range("O92").Value = 15000
$ws.Range("P92").Value = 14500
$ws.Range("Q92").Value = '$/caja 14 unidades'
$ws.Range("S92").Value = 1036
$ws.Range("T92").Value = 14

$ws.Range("D93").Value = 44222
$ws.Range("L93").Value = 'Primera'
$ws.Range("M93").Value = 50
$ws.Range("N93").Value = 15000
$ws.Range("O93").Value = 16000
$ws.Range("P93").Value = 15400
$ws.Range("Q93").Value = '$/caja 12 unidades'
$ws.Range("S93").Value = 1283
$ws.Range("T93").Value = 12

$ws.Range("D94").Value = 44398
$ws.Range("L94").Value = 'Segunda'
$ws.Range("M94").Value = 80
$ws.Range("N94").Value = 16000
$ws.Range("O94").Value = 16000
$ws.Range("P94").Value = 16000
$ws.Range("Q94").Value = '$/caja 14 unidades'
$ws.Range("S94").Value = 1143
$ws.Range("T94").Value = 14

$ws.Range("D95").Value = 44264
$ws.Range("L95").Value = 'Segunda'
$ws.Range("M95").Value = 60
$ws.Range("N95").Value = 16000
$ws.Range("O95").Value = 17000
$ws.Range("P95").Value = 16500
$ws.Range("Q95").Value = '$/caja 14 unidades'
$ws.Range("S95").Value = 1179
$ws.Range("T95").Value = 14

$ws.Range("D96").Value = 44210
$ws.Range("L96").Value = 'Primera'
$ws.Range("M96").Value = 60
$ws.Range("N96").Value = 13000
$ws.Range("O96").Value = 14000
$ws.Range("P96").Value = 13333
$ws.Range("Q96").Value = '$/caja 12 unidades'
$ws.Range("S96").Value = 1111
$ws.Range("T96").Value = 12

$ws.Range("D97").Value = 44253
$ws.Range("L97").Value = 'Primera'
$ws.Range("M97").Value = 50
$ws.Range("N97").Value = 14500
$ws.Range("O97").Value = 15000
$ws.Range("P97").Value = 14800
$ws.Range("Q97").Value = '$/caja 12 unidades'
$ws.Range("S97").Value = 1233
$ws.Range("T97").Value = 12

$ws.Range("D98").Value = 44341
$ws.Range("L98").Value = 'Segunda'
$ws.Range("M98").Value = 120
$ws.Range("N98").Value = 15500
$ws.Range("O98").Value = 16000
$ws.Range("P98").Value = 15750
$ws.Range("Q98").Value = '$/caja 14 unidades'
$ws.Range("S98").Value = 1125
$ws.Range("T98").Value = 14

$ws.Range("D99").Value = 44414
$ws.Range("L99").Value = 'Primera'
$ws.Range("M99").Value = 100
$ws.Range("N99").Value = 18000
$ws.Range("O99").Value = 19000
$ws.Range("P99").Value = 18500
$ws.Range("Q99").Value = '$/caja 12 unidades'
$ws.Range("S99").Value = 1542
$ws.Range("T99").Value = 12

$ws.Range("D100").Value = 44414
$ws.Range("L100").Value = 'Segunda'
$ws.Range("M100").Value = 100
$ws.Range("N100").Value = 18000
$ws.Range("O100").Value = 19000
$ws.Range("P100").Value = 18500
$ws.Range("Q100").Value = '$/caja 14 unidades'
$ws.Range("S100").Value = 1321
$ws.Range("T100").Value = 14

$ws.Range("D101").Value = 44414
$ws.Range("L101").Value = 'Tercera'
$ws.Range("M101").Value = 40
$ws.Range("N101").Value = 18000
$ws.Range("O101").Value = 18000
$ws.Range("P101").Value = 18000
$ws.Range("Q101").Value = '$/caja 16 unidades'
$ws.Range("S101").Value = 1125
$ws.Range("T101").Value = 16

$ws.Range("D102").Value = 44301
$ws.Range("L102").Value = 'Segunda'
$ws.Range("M102").Value = 120
$ws.Range("N102").Value = 15000
$ws.Range("O102").Value = 16000
$ws.Range("P102").Value = 15500
$ws.Range("Q102").Value = '$/caja 14 unidades'
$ws.Range("S102").Value = 1107
$ws.Range("T102").Value = 14

$ws.Range("D103").Value = 44273
$ws.Range("L103").Value = 'Segunda'
$ws.Range("M103").Value = 120
$ws.Range("N103").Value = 15000
$ws.Range("O103").Value = 16000
$ws.Range("P103").Value = 15500
$ws.Range("Q103").Value = '$/caja 14 unidades'
$ws.Range("S103").Value = 1107
$ws.Range("T103").Value = 14

$ws.Range("D104").Value = 44342
$ws.Range("L104").Value = 'Segunda'
$ws.Range("M104").Value = 120
$ws.Range("N104").Value = 15000
$ws.Range("O104").Value = 16000
$ws.Range("P104").Value = 15500
$ws.Range("Q104").Value = '$/caja 14 unidades'
$ws.Range("S104").Value = 1107
$ws.Range("T104").Value = 14

$ws.Range("D105").Value = 44294
$ws.Range("L105").Value = 'Segunda'
$ws.Range("M105").Value = 40
$ws.Range("N105").Value = 15000
$ws.Range("O105").Value = 16000
$ws.Range("P105").Value = 15500
$ws.Range("Q105").Value = '$/caja 14 unidades'
$ws.Range("S105").Value = 1107
$ws.Range("T105").Value = 14

$ws.Range("D106").Value = 44379
$ws.Range("L106").Value = 'Segunda'
$ws.Range("M106").Value = 120
$ws.Range("N106").Value = 16000
$ws.Range("O106").Value = 17000
$ws.Range("P106").Value = 16500
$ws.Range("Q106").Value = '$/caja 14 unidades'
$ws.Range("S106").Value = 1179
$ws.Range("T106").Value = 14

$ws.Range("D107").Value = 44302
$ws.Range("L107").Value = 'Segunda'
$ws.Range("M107").Value = 60
$ws.Range("N107").Value = 15000
$ws.Range("O107").Value = 16000
$ws.Range("P107").Value = 15500
$ws.Range("Q107").Value = '$/caja 14 unidades'
$ws.Range("S107").Value = 1107
$ws.Range("T107").Value = 14

$ws.Range("D108").Value = 44315
$ws.Range("L108").Value = 'Segunda'
$ws.Range("M108").Value = 120
$ws.Range("N108").Value = 15000
$ws.Range("O108").Value = 16000
$ws.Range("P108").Value = 15500
$ws.Range("Q108").Value = '$/caja 14 unidades'
$ws.Range("S108").Value = 1107
$ws.Range("T108").Value = 14

$ws.Range("D109").Value = 44446
$ws.Range("L109").Value = 'Primera'
$ws.Range("M109").Value = 100
$ws.Range("N109").Value = 19000
$ws.Range("O109").Value = 20000
$ws.Range("P109").Value = 19500
$ws.Range("Q109").Value = '$/caja 12 unidades'
$ws.Range("S109").Value = 1625
$ws.Range("T109").Value = 12

$ws.Range("D110").Value = 44446
$ws.Range("L110").Value = 'Segunda'
$ws.Range("M110").Value = 60
$ws.Range("N110").Value = 19000
$ws.Range("O110").Value = 20000
$ws.Range("P110").Value = 19500
$ws.Range("Q110").Value = '$/caja 14 unidades'
$ws.Range("S110").Value = 1393
$ws.Range("T110").Value = 14

$ws.Range("D111").Value = 44313
$ws.Range("L111").Value = 'Segunda'
$ws.Range("M111").Value = 60
$ws.Range("N111").Value = 15000
$ws.Range("O111").Value = 16000
$ws.Range("P111").Value = 15500
$ws.Range("Q111").Value = '$/caja 14 unidades'
$ws.Range("S111").Value = 1107
$ws.Range("T111").Value = 14

$ws.Range("D112").Value = 44329
$ws.Range("L112").Value = 'Segunda'
$ws.Range("M112").Value = 120
$ws.Range("N112").Value = 16000
$ws.Range("O112").Value = 17000
$ws.Range("P112").Value = 16500
$ws.Range("Q112").Value = '$/caja 14 unidades'
$ws.Range("S112").Value = 1179
$ws.Range("T112").Value = 14

$ws.Range("D113").Value = 44161
$ws.Range("L113").Value = 'Segunda'
$ws.Range("M113").Value = 30
$ws.Range("N113").Value = 21000
$ws.Range("O113").Value = 22000
$ws.Range("P113").Value = 21500
$ws.Range("Q113").Value = '$/caja 14 unidades'
$ws.Range("S113").Value = 1536
$ws.Range("T113").Value = 14

$ws.Range("D114").Value = 44251
$ws.Range("L114").Value = 'Primera'
$ws.Range("M114").Value = 100
$ws.Range("N114").Value = 14000
$ws.Range("O114").Value = 15000
$ws.Range("P114").Value = 14650
$ws.Range("Q114").Value = '$/caja 12 unidades'
$ws.Range("S114").Value = 1221
$ws.Range("T114").Value = 12

$ws.Range("D115").Value = 44319
$ws.Range("L115").Value = 'Segunda'
$ws.Range("M115").Value = 60
$ws.Range("N115").Value = 16000
$ws.Range("O115").Value = 17000
$ws.Range("P115").Value = 16500
$ws.Range("Q115").Value = '$/caja 14 unidades'
$ws.Range("S115").Value = 1179
$ws.Range("T115").Value = 14

$ws.Range("D116").Value = 44175
$ws.Range("L116").Value = 'Primera'
$ws.Range("M116").Value = 55
$ws.Range("N116").Value = 19000
$ws.Range("O116").Value = 20000
$ws.Range("P116").Value = 19545
$ws.Range("Q116").Value = '$/caja 12 unidades'
$ws.Range("S116").Value = 1629
$ws.Range("T116").Value = 12

$ws.Range("D117").Value = 44175
$ws.Range("L117").Value = 'Segunda'
$ws.Range("M117").Value = 40
$ws.Range("N117").Value = 19500
$ws.Range("O117").Value = 20000
$ws.Range("P117").Value = 19750
$ws.Range("Q117").Value = '$/caja 14 unidades'
$ws.Range("S117").Value = 1411
$ws.Range("T117").Value = 14

$ws.Range("D118").Value = 44376
$ws.Range("L118").Value = 'Segunda'
$ws.Range("M118").Value = 120
$ws.Range("N118").Value = 16000
$ws.Range("O118").Value = 17000
$ws.Range("P118").Value = 16500
$ws.Range("Q118").Value = '$/caja 14 unidades'
$ws.Range("S118").Value = 1179
$ws.Range("T118").Value = 14

$ws.Range("D119").Value = 44279
$ws.Range("L119").Value = 'Segunda'
$ws.Range("M119").Value = 120
$ws.Range("N119").Value = 15000
$ws.Range("O119").Value = 16000
$ws.Range("P119").Value = 15500
$ws.Range("Q119").Value = '$/caja 14 unidades'
$ws.Range("S119").Value = 1107
$ws.Range("T119").Value = 14

$ws.Range("D120").Value = 44223
$ws.Range("L120").Value = 'Segunda'
$ws.Range("M120").Value = 60
$ws.Range("N120").Value = 14000
$ws.Range("O120").Value = 15000
$ws.Range("P120").Value = 14500
$ws.Range("Q120").Value = '$/caja 14 unidades'
$ws.Range("S120").Value = 1036
$ws.Range("T120").Value = 14

$ws.Range("D121").Value = 44314
$ws.Range("L121").Value = 'Segunda'
$ws.Range("M121").Value = 120
$ws.Range("N121").Value = 15000
$ws.Range("O121").Value = 16000
$ws.Range("P121").Value = 15500
$ws.Range("Q121").Value = '$/caja 14 unidades'
$ws.Range("S121").Value = 1107
$ws.Range("T121").Value = 14

$ws.Range("D122").Value = 44399
$ws.Range("L122").Value = 'Segunda'
$ws.Range("M122").Value = 120
$ws.Range("N122").Value = 17000
$ws.Range("O122").Value = 18000
$ws.Range("P122").Value = 17500
$ws.Range("Q122").Value = '$/caja 14 unidades'
$ws.Range("S122").Value = 1250
$ws.Range("T122").Value = 14

